$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.514.01'
$ws.Range("E2").Value = '  -0.04%  '

$ws.Range("D3").Value = '2.644.58'
$ws.Range("E3").Value = '  +0.10%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '602.55'
$ws.Range("E5").Value = '  +2.04%  '

$ws.Range("D6").Value = '146.97'
$ws.Range("E6").Value = '  +1.12%  '

$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("D8").Value = '0.587'
$ws.Range("E8").Value = '  -0.42%  '

$ws.Range("E9").Value = '  +1.07%  '

$ws.Range("E10").Value = '  -0.98%  '

$ws.Range("D11").Value = '0.368'
$ws.Range("E11").Value = '  +4.19%  '

$ws.Range("E12").Value = '  +0.01%  '

$ws.Range("D13").Value = '27.36'
$ws.Range("E13").Value = '  -0.74%  '

$ws.Range("D14").Value = '3.120.17'
$ws.Range("E14").Value = '  +0.22%  '

$ws.Range("D15").Value = '63.349.48'
$ws.Range("E15").Value = '  -0.16%  '

$ws.Range("E16").Value = '  -0.34%  '

$ws.Range("D17").Value = '2.613.98'
$ws.Range("E17").Value = '  -3.05%  '

$ws.Range("D18").Value = '11.46'
$ws.Range("E18").Value = '  +1.41%  '

$ws.Range("D19").Value = '4.53'
$ws.Range("E19").Value = '  +4.02%  '

$ws.Range("D20").Value = '340.78'
$ws.Range("E20").Value = '  -0.10%  '

$ws.Range("D21").Value = '6.95'
$ws.Range("E21").Value = '  +3.67%  '

$ws.Range("E22").Value = '  -0.01%  '

$ws.Range("D23").Value = '5.57'
$ws.Range("E23").Value = '  -3.22%  '

$ws.Range("D24").Value = '66.70'
$ws.Range("E24").Value = '  -1.61%  '

$ws.Range("E25").Value = '  +0.44%  '

$ws.Range("E26").Value = '  +4.74%  '

$ws.Range("E27").Value = '  -2.02%  '

$ws.Range("D28").Value = '0.164'
$ws.Range("E28").Value = '  -1.53%  '

$ws.Range("D29").Value = '547.32'
$ws.Range("E29").Value = '  -1.00%  '

$ws.Range("E30").Value = '  -0.08%  '

$ws.Range("E31").Value = '  +0.89%  '

$ws.Range("E32").Value = '  +4.20%  '

$ws.Range("E33").Value = '  -2.99%  '

$ws.Range("D34").Value = '0.0₃0807'
$ws.Range("E34").Value = '  +0.05%  '

$ws.Range("D35").Value = '5.26'
$ws.Range("E35").Value = '  +8.01%  '

$ws.Range("D36").Value = '167.92'
$ws.Range("E36").Value = '  -4.07%  '

$ws.Range("E37").Value = '  +1.85%  '

$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  -0.04%  '

$ws.Range("D39").Value = '19.08'
$ws.Range("E39").Value = '  -0.08%  '

$ws.Range("D40").Value = '1.91'
$ws.Range("E40").Value = '  +7.15%  '

$ws.Range("E41").Value = '  -0.06%  '

$ws.Range("D42").Value = '169.52'
$ws.Range("E42").Value = '  -0.17%  '

$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").Value = '22.73'
$ws.Range("E43").Value = '  +2.87%  '

$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").Value = '3.78'
$ws.Range("E44").Value = '  +1.30%  '

$ws.Range("E45").Value = '  +4.97%  '

$ws.Range("D46").Value = '0.626'
$ws.Range("E46").Value = '  -0.20%  '

$ws.Range("D47").Value = '0.0247'
$ws.Range("E47").Value = '  +3.39%  '

$ws.Range("D48").Value = '0.0963'
$ws.Range("E48").Value = '  +0.42%  '

$ws.Range("E49").Value = '  +0.68%  '

$ws.Range("D50").Value = '1.80'
$ws.Range("E50").Value = '  +5.70%  '

$ws.Range("E51").Value = '  -0.80%  '
